$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Components" (sheet 1)
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Components")

# Fix wording of the humidity sensor note (row 4)
$ws.Range("G4").Value = "rel humidity & ext air temp - ext mount"

# Insert a new component row for the voltage-reg bypass capacitor (new row 10,
# pushing the existing rows 10-15 down to 11-16)
$ws.Rows.Item(10).Insert()

$ws.Range("A10").Value = "10nf ceramic"
$ws.Range("B10").Value = 1
$ws.Range("C10").Value = 0.122
$ws.Range("D10").Formula = "=B10*C10"
$ws.Range("E10").Value = "Farnell"
$ws.Range("F10").Value = "0805"
$ws.Range("G10").Value = "Voltage reg bypass cap. C0G/NP0 as specified in app notes"
$ws.Range("H10").Value = "http://uk.farnell.com/avx/0805ya103jat2a/cap-mlcc-c0g-np0-10nf-16v-0805/dp/2332815"

# Tidy up package notes that were previously "N/A" / missing for a few parts
$ws.Range("F13").Value = "O/B"
$ws.Range("F15").Value = "O/B"
$ws.Range("F16").Value = "O/B"

# Make room for two more new component rows just below "Motors" (row 16),
# above the blank separator row that precedes the Subtotal row
$ws.Range("17:18").Insert()

$ws.Range("A17").Value = "10K resistors"
$ws.Range("B17").Value = 2
$ws.Range("C17").Value = 0.014
$ws.Range("D17").Formula = "=B17*C17"
$ws.Range("E17").Value = "Farnell"
$ws.Range("F17").Value = "0805"
$ws.Range("G17").Value = "I2C SDA&SCL pullups"
$ws.Range("H17").Value = "http://uk.farnell.com/te-connectivity/crgh0805f10k/resistor-power-10k-0-33w-1-0805/dp/2332084RL"

$ws.Range("A18").Value = "Diode"
$ws.Range("B18").Value = 1
$ws.Range("C18").Value = 0.078
$ws.Range("D18").Formula = "=B18*C18"
$ws.Range("E18").Value = "Farnell"
$ws.Range("F18").Value = "SOD-323"
$ws.Range("G18").Value = "Boost power supply diode. 1A forward current. "
$ws.Range("H18").Value = "http://uk.farnell.com/nxp/bat760/diode-schottky-sod-323/dp/8734593"

# Extend the Subtotal formula to cover the two freshly added rows
$ws.Range("C20").Formula = "=SUM(D2:D18)"

# ---------------------------------------------------------------------------
# Sheet "PCB names and values" (sheet 2)
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("PCB names and values")

$ws2.Range("A6").Value = "C1"
$ws2.Range("B6").Value = "Capacitor for U2"
$ws2.Range("C6").Value = "4.7uF"

$ws2.Range("A7").Value = "C2"
$ws2.Range("B7").Value = "Capacitor for U2"
$ws2.Range("C7").Value = "4.7uF"

$ws2.Range("A8").Value = "C3"
$ws2.Range("B8").Value = "Capacitor for U2"
$ws2.Range("C8").Value = "10nf"

$ws2.Range("C8").Select()

# ---------------------------------------------------------------------------
# Active sheet / selection bookkeeping: the "Components" sheet is now the
# active tab (was "PCB names and values" before), with G18 selected.
# ---------------------------------------------------------------------------
$ws.Activate()
$ws.Range("G18").Select()
